$d = $word.ActiveDocument

$sec = $d.Sections.First
$footer = $sec.Footers.Item(1)
$fr = $footer.Range

# "App Name – " -> "UC Schedule – "  (drop "Name", keep the en dash)
$fr.Find.Execute("App Name " + [char]0x2013 + " ", $true, $false, $false, $false, $false,
                  $true, 1, $false, "UC Schedule " + [char]0x2013 + " ", 2)

# Refresh the cached PAGE field result (2 -> 1) without disturbing the
# surrounding text runs.
$pageField = $footer.Range.Fields.Item(1)
$pageField.Result.Find.Execute("2", $true, $false, $false, $false, $false,
                                $true, 1, $false, "1", 2)

Write-Output $footer.Range.Text
